$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ui")

# The "ui" sheet lists HP-recovery related messages; append two new detail rows
# (29, 30) right after the existing table, matching the formatting/formula
# pattern used by the other rows (3-28): column A is the running row number
# ("=ROW()-2"), column B is the message text.

# Copy formatting (font/fill/border/number format) and row height from the
# last existing row (28) down into the two new rows.
$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Rows.Item(29).RowHeight = $ws.Rows.Item(28).RowHeight

$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A30:B30").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Rows.Item(30).RowHeight = $ws.Rows.Item(28).RowHeight

# Row-number formula, continuing the existing series.
$ws.Range("A29:A30").Formula = "=ROW()-2"

# New message descriptions.
$ws.Range("B29").Value = "HPを10%回復します"
$ws.Range("B30").Value = "HPを最大値まで回復します"
